# Split the old "C1 (13:00-15:00)" slot on row 9 into two slots:
#   row 9  -> "S1 (07:00-09:00)"  keeps the "Tieng Anh chuyen nganh" class,
#             now placed under "Thu 3" (column E) instead of "Thu 5" (column G)
#   row 10 -> new row "C2 (15:00-17:00)" holding the "Ky nang mem" class,
#             now placed under "Thu 6" (column H) instead of "Thu 4" (column F)
# Applied identically across every weekly sheet (Tuan_1 .. Tuan_15).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # 1) Clone row 9 (values + formatting) into the new row 10.
    $ws.Range("A9:I9").Copy($ws.Range("A10:I10"))
    $ws.Rows.Item(10).RowHeight = $ws.Rows.Item(9).RowHeight

    # 2) On row 10, move the "Ky nang mem" content from F10 to H10 (Thu 6 / C2),
    #    then blank out F10/G10 back to the plain (unhighlighted) look.
    $ws.Range("F10").Copy($ws.Range("H10"))
    $ws.Range("D10").Copy()
    $ws.Range("F10").PasteSpecial($xlPasteFormats)
    $ws.Range("G10").PasteSpecial($xlPasteFormats)
    $ws.Range("F10").ClearContents()
    $ws.Range("G10").ClearContents()
    $excel.CutCopyMode = $false

    $ws.Range("A10").Value = "C2" + [char]10 + "(15:00-17:00)"

    # 3) On row 9, move the "Tieng Anh chuyen nganh" content from G9 to E9
    #    (Thu 3 / S1), then blank out F9/G9 back to the plain look.
    $ws.Range("G9").Copy($ws.Range("E9"))
    $ws.Range("D9").Copy()
    $ws.Range("F9").PasteSpecial($xlPasteFormats)
    $ws.Range("G9").PasteSpecial($xlPasteFormats)
    $ws.Range("F9").ClearContents()
    $ws.Range("G9").ClearContents()
    $excel.CutCopyMode = $false

    $ws.Range("A9").Value = "S1" + [char]10 + "(07:00-09:00)"
}
